# Fix the 2050 column-header label (it was holding a stray numeric value
# left over from a bad paste/formula) and drop the "Total" row that each
# table carried.
#
# Helper: write $text into $addr as genuine text (not a number), matching
# how the other header cells in the same row are already stored, and
# without disturbing the destination cell's existing style. We stage the
# text in a scratch cell far outside the used range, force it to Text via
# NumberFormat, copy it, and paste *values only* into the destination so
# the destination keeps its own formatting (PasteSpecial values-only does
# not carry the scratch cell's number format along). The scratch cell is
# then fully cleared so it leaves no trace in the sheet's used range.
function Set-TextValue($ws, $addr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook

# Sheets whose E1 header should read "2050" and whose trailing "Total" row
# (row 13) should be removed entirely.
$sheetsWith2050 = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)
foreach ($name in $sheetsWith2050) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextValue $ws "E1" "2050"
    $ws.Rows.Item(13).Delete()
}

# "Potencia Incremental - SIN(MW)" uses period ranges as headers, so its
# E1 should read "2041-2050" instead of a bare year. Same Total-row removal.
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextValue $ws4 "E1" "2041-2050"
$ws4.Rows.Item(13).Delete()

# "Emissoes Totais (MtCO2eq)" only needs the header fix -- it never had a
# Total row to begin with.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextValue $ws5 "E1" "2050"

# "Custo Total (bilhões de R$)" has no mis-typed header, just its own
# Total row (row 4) to drop.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
